$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "serorreversion" -> "seroreversion" in the Spanish header row (row 9)
$ws.Range("F9").Value = "Tasa de seroreversion"
$ws.Range("G9").Value = "Tasa de seroreversión Rhat"

# Update the selection to match the second (Spanish) table block
$ws.Range("A9:G13").Select()
